# Revert "Revert "Add retry in case of server error""
# Re-adds the ServerErrorFailure localization row (row 43) to the
# "Localization" sheet, growing Table13 by one row, and leaves that
# sheet as the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Insert a new row before the existing row 43, shifting rows 43-108 down
# by one (row 43 "ConfirmNumerousRequests" becomes row 44, etc.)
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row with the restored strings.
$ws.Cells.Item(43, 1).Value2 = "ServerErrorFailure"
$ws.Cells.Item(43, 2).Value2 = "HTTP Request failed due to server error issues."
$ws.Cells.Item(43, 3).Value2 = "サーバーエラーの問題のため、リクエストが失敗しました。"

# Match the wrap-text formatting used by the other Name/EN/JA rows
# (column A stays unwrapped, columns B/C wrap like every other entry).
$ws.Cells.Item(43, 2).WrapText = $true
$ws.Cells.Item(43, 3).WrapText = $true

# Grow the worksheet table (Table13) so it covers the extra row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C109"))

# The "Localization" sheet becomes the active tab.
$ws.Activate()
